# Apply the published-CDA-FHIR-logical-model update (patches #241) to the
# "StructureDefinition-Entity" workbook.
#
# Sheet "Metadata" holds Property/Value pairs (A/B columns).
# Sheet "Elements" holds the element table (row 1 = headers).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: "No display for ContactDetail" -> full HL7 contact string
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# Entity.classCode row (row 12): Min / Base Min go from 1 to 0
$elements.Range("F12").Value = "0"
$elements.Range("AG12").Value = "0"

# Entity.code row (row 15): Binding Value Set URL update
$elements.Range("Z15").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAEntityCode"
